# Update the prompts file: insert a \QUESTIONBREAK marker right before the
# repeated "\item According to the PDF pages ..." instruction in every
# generated-prompt cell (column B), then leave the selection on B1 and
# nudge column B's formatting (long text column -> wrap) which is what
# produced the extra cell-style slot in the saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.UsedRange
$rng.Replace("Use this exact format: \item", "Use this exact format: \QUESTIONBREAK \item")

$colB = $ws.Range("B1:B20")
$colB.WrapText = $true

$ws.Range("B1").Select()
